$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 2105981.08
$ws.Range("C7").Value = -53.38101999025433
$ws.Range("D7").Value = 1989
$ws.Range("E7").Value = 1989
$ws.Range("F7").Value = 1058.814017094017
$ws.Range("G7").Value = 9.269826448182172
